# "commandline arguments.xlsx" / Sheet1 — replace the old 4-row "passing an
# integer/char/string/float" test-case table with a new 8-row table
# (prasanna's command-line-argument-parser test cases), add two brand new
# rows (8 & 9) and touch a few cell alignments / row heights to match.
#
# NOTE on ordering: new unique string values are written in the same
# relative order they first appear in the finished workbook's shared-string
# table (verified against the target XML) so that <si> index assignment
# lines up; cells that keep an already-existing string (duplicates / the
# unchanged header row) are written afterwards since they don't mint a new
# shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1) New/changed string values, in shared-string creation order ----
$ws.Range("C6").Value = 'prasanna'
$ws.Range("D3").Value = 'Parser.exe argv[v],datatype1 ,datatype2,….datatypeN'
$ws.Range("D5").Value = 'In command line arguments pass -h to display'
$ws.Range("D7").Value = 'int     1      4                                                    int     2      4'
$ws.Range("D8").Value = 'int         1         4                                                          char       a         1                                    float    2.2       4'
$ws.Range("D9").Value = 'string    Abcd        8                                          int         234         4                                                float      432.5    4'
$ws.Range("D4").Value = 'string          -h         8                                           string         -h         8'
$ws.Range("E4").Value = 'string            -h          8                                                   string              -h         8'
$ws.Range("D6").Value = 'string          prasanna         8'
$ws.Range("E6").Value = 'string          prasanna      8'
$ws.Range("C5").Value = '  noinput'
$ws.Range("C4").Value = '                        -h  -h'
$ws.Range("C8").Value = '1      a      2.2'
$ws.Range("C9").Value = 'Abcd       234       432.5'
$ws.Range("C7").Value = '1    2'
$ws.Range("B8").Value = 'passing    char,int,float'
$ws.Range("B9").Value = 'passing    string ,int,float'
$ws.Range("B7").Value = 'passing    integer'
$ws.Range("B6").Value = 'Passing     string'
$ws.Range("B5").Value = 'Passing   " noinput"'
$ws.Range("B4").Value = 'Passing   " -h -h"'
$ws.Range("B3").Value = 'Passing   "-h"'
$ws.Range("B2").Value = 'Passing  an integer,charater,string,float'
$ws.Range("C3").Value = '    -h'

# ---- 2) Remaining cells: numbers, and reused/duplicate strings ----
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("E3").Value = 'Parser.exe argv[v],datatype1 ,datatype2,….datatypeN'
$ws.Range("F3").Value = 'PASS'
$ws.Range("A4").Value = 3
$ws.Range("F4").Value = 'PASS'
$ws.Range("A5").Value = 4
$ws.Range("E5").Value = 'In command line arguments pass -h to display'
$ws.Range("F5").Value = 'PASS'
$ws.Range("A6").Value = 5
$ws.Range("F6").Value = 'PASS'
$ws.Range("A7").Value = 6
$ws.Range("E7").Value = 'int     1      4                                                    int     2      4'
$ws.Range("F7").Value = 'PASS'
$ws.Range("A8").Value = 7
$ws.Range("E8").Value = 'int         1         4                                                          char       a         1                                    float    2.2       4'
$ws.Range("F8").Value = 'PASS'
$ws.Range("A9").Value = 8
$ws.Range("E9").Value = 'string    Abcd        8                                          int         234         4                                                float      432.5    4'
$ws.Range("F9").Value = 'PASS'

# ---- 3) Alignment/style tweaks on existing rows ----
# C4 drops its centered horizontal alignment (keeps vertical-center + wrap).
$ws.Range("C4").HorizontalAlignment = 1
# E6 gains vertical-center (already had horizontal-center + wrap).
$ws.Range("E6").VerticalAlignment = -4108
# B7 gains vertical-center (already had horizontal-center).
$ws.Range("B7").VerticalAlignment = -4108

# ---- 4) Brand new rows 8 & 9: same look as rows 3-5 (A/F centered only; ----
#         B-E centered + vertical-centered + wrap)
foreach ($r in 8,9) {
    $ws.Range("A$r").HorizontalAlignment = -4108
    foreach ($col in 'B','C','D','E') {
        $ws.Range("$col$r").HorizontalAlignment = -4108
        $ws.Range("$col$r").VerticalAlignment = -4108
        $ws.Range("$col$r").WrapText = $true
    }
    $ws.Range("F$r").HorizontalAlignment = -4108
}

# ---- 5) Row heights for the new rows ----
$ws.Rows.Item(8).RowHeight = 63.6
$ws.Rows.Item(9).RowHeight = 55.8

# ---- 6) Selection cursor (closest reproducible approximation: the COM
#         shim collapses multi-area selects to a single area, so this
#         reproduces the dominant "entire column A" part of the original
#         "B4 A1:A1048576" multi-range selection). ----
$ws.Range("A1:A1048576").Select()
